# Update the "cryptos" price/volume table with the latest scraped values.
# Note: price column (D) values that look like plain numbers (e.g. "10.99")
# are prefixed with a leading apostrophe so Excel stores them as text
# instead of silently reinterpreting them as numeric values (matching the
# original inlineStr/text cell content).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.182.10"
$ws.Range("E2").Value = "  +0.25%  "
$ws.Range("D3").Value = "1.786.12"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'225.87"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("E6").Value = "  +0.39%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").Value = "'32.10"
$ws.Range("E8").Value = "  -0.31%  "
$ws.Range("D9").Value = "'0.292"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "'0.0686"
$ws.Range("E10").Value = "  -0.12%  "
$ws.Range("E11").Value = "  +1.24%  "
$ws.Range("D12").Value = "2.044.34"
$ws.Range("E12").Value = "  -0.13%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.789.14"
$ws.Range("E13").Value = "  -0.16%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'10.99"
$ws.Range("E14").Value = "  -3.23%  "
$ws.Range("D15").Value = "'0.625"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "34.156.30"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "'67.78"
$ws.Range("E18").Value = "  -0.43%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "'245.26"
$ws.Range("E19").Value = "  +0.73%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0796"
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").Value = "'10.99"
$ws.Range("E21").Value = "  +0.76%  "
$ws.Range("E22").Value = "  +0.25%  "
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("E24").Value = "  -0.42%  "
$ws.Range("D25").Value = "'161.26"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("D28").Value = "'0.114"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").Value = "'1.00"
$ws.Range("E29").Value = "  +0.30%  "
$ws.Range("B30").Value = "PancakeSwap"
$ws.Range("C30").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D30").Value = "'1.23"
$ws.Range("E30").Value = "  -0.51%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").Value = "'0.0518"
$ws.Range("E31").Value = "  -0.48%  "
$ws.Range("E32").Value = "  +2.60%  "
$ws.Range("D33").Value = "'3.74"
$ws.Range("E34").Value = "  -2.03%  "
$ws.Range("D35").Value = "1.444.38"
$ws.Range("E35").Value = "  +3.00%  "
$ws.Range("D36").Value = "'2.61"
$ws.Range("E36").Value = "  +11.25%  "
$ws.Range("D37").Value = "'0.657"
$ws.Range("E37").Value = "  +1.06%  "
$ws.Range("E38").Value = "  +0.75%  "
$ws.Range("D39").Value = "'1.05"
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("E41").Value = "  +1.68%  "
$ws.Range("D42").Value = "'13.96"
$ws.Range("E42").Value = "  +4.55%  "
$ws.Range("E43").Value = "  +0.72%  "
$ws.Range("D44").Value = "'0.917"
$ws.Range("E44").Value = "  -0.34%  "
$ws.Range("D45").Value = "'0.0515"
$ws.Range("E45").Value = "  +1.15%  "
$ws.Range("E46").Value = "  +0.82%  "
$ws.Range("E47").Value = "  +0.51%  "
$ws.Range("D48").Value = "1.943.41"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D49").Value = "'104.93"
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("E51").Value = "  -5.80%  "
